$wb = $excel.ActiveWorkbook

# --- Sheet: ALC ---
$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H40").Value = 1833.3334
$ws.Range("I40").Value = 1833.3334
$ws.Range("J40").Value = 0
$ws.Range("K40").Value = 1833.3334
$ws.Range("L40").Value = 0
$ws.Range("M40").Value = -1658.3334
$ws.Range("N40").Value = ""
$ws.Range("H70").Value = 3856.25
$ws.Range("J70").Value = 4097.1113
$ws.Range("L70").Value = 12291.3339
$ws.Range("N70").Value = -12831.3339
$ws.Range("H73").Value = 3856.25
$ws.Range("J73").Value = 4097.1113
$ws.Range("L73").Value = 12291.3339
$ws.Range("N73").Value = -14163.3339
$ws.Range("H125").Value = 3640.1428
$ws.Range("I125").Value = 1912.625
$ws.Range("J125").Value = 5943.5
$ws.Range("K125").Value = 17213.625
$ws.Range("L125").Value = 53491.5
$ws.Range("M125").Value = -14753.625
$ws.Range("N125").Value = -58411.5
$ws.Range("H127").Value = 1474.25
$ws.Range("I127").Value = 965.6667
$ws.Range("K127").Value = 2897.0001
$ws.Range("M127").Value = 2062.9999
$ws.Range("H137").Value = 3158
$ws.Range("I137").Value = 3073.1
$ws.Range("J137").Value = 4007
$ws.Range("K137").Value = 9219.299999999999
$ws.Range("L137").Value = 12021
$ws.Range("M137").Value = -6669.299999999999
$ws.Range("N137").Value = -17121

# --- Sheet: ARM ---
$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H3").Value = 1222
$ws.Range("J3").Value = 3656
$ws.Range("L3").Value = 3656
$ws.Range("N3").Value = -3886
$ws.Range("H16").Value = 17857.834
$ws.Range("I16").Value = 50100
$ws.Range("J16").Value = 1736.75
$ws.Range("K16").Value = 50100
$ws.Range("L16").Value = 1736.75
$ws.Range("M16").Value = -49813
$ws.Range("N16").Value = -2310.75
$ws.Range("H27").Value = 0
$ws.Range("J27").Value = 0
$ws.Range("L27").Value = 0
$ws.Range("N27").Value = ""
$ws.Range("H30").Value = 2050
$ws.Range("I30").Value = 2333.3333
$ws.Range("K30").Value = 2333.3333
$ws.Range("M30").Value = -2183.3333
$ws.Range("H61").Value = 2548.6667
$ws.Range("I61").Value = 2458.4
$ws.Range("K61").Value = 2458.4
$ws.Range("M61").Value = -2246.4
$ws.Range("H69").Value = 250000
$ws.Range("J69").Value = 250000
$ws.Range("L69").Value = 250000
$ws.Range("N69").Value = -251498
$ws.Range("H72").Value = 250000
$ws.Range("J72").Value = 250000
$ws.Range("L72").Value = 750000
$ws.Range("N72").Value = -757488
$ws.Range("H74").Value = 1270.875
$ws.Range("I74").Value = 1293.2858
$ws.Range("J74").Value = 1114
$ws.Range("K74").Value = 1293.2858
$ws.Range("L74").Value = 1114
$ws.Range("M74").Value = -419.2858000000001
$ws.Range("N74").Value = -2862
$ws.Range("H77").Value = 1270.875
$ws.Range("I77").Value = 1293.2858
$ws.Range("J77").Value = 1114
$ws.Range("K77").Value = 6466.429
$ws.Range("L77").Value = 5570
$ws.Range("M77").Value = -2098.429
$ws.Range("N77").Value = -14306
$ws.Range("H132").Value = 2058
$ws.Range("I132").Value = 2058
$ws.Range("K132").Value = 6174
$ws.Range("M132").Value = -3644
$ws.Range("H136").Value = 2548.6667
$ws.Range("I136").Value = 2458.4
$ws.Range("K136").Value = 7375.200000000001
$ws.Range("M136").Value = -4825.200000000001

# --- Sheet: BSM ---
$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 4462.5884
$ws.Range("I86").Value = 3897.0833
$ws.Range("J86").Value = 5819.8
$ws.Range("K86").Value = 3897.0833
$ws.Range("L86").Value = 5819.8
$ws.Range("M86").Value = -2774.0833
$ws.Range("N86").Value = -8065.8
$ws.Range("H89").Value = 4462.5884
$ws.Range("I89").Value = 3897.0833
$ws.Range("J89").Value = 5819.8
$ws.Range("K89").Value = 19485.4165
$ws.Range("L89").Value = 29099
$ws.Range("M89").Value = -13869.4165
$ws.Range("N89").Value = -40331
$ws.Range("H94").Value = 1752.9
$ws.Range("I94").Value = 1752.9
$ws.Range("J94").Value = 0
$ws.Range("K94").Value = 1752.9
$ws.Range("L94").Value = 0
$ws.Range("M94").Value = -1301.9
$ws.Range("N94").Value = ""
$ws.Range("H105").Value = 3978.125
$ws.Range("I105").Value = 3832.2856
$ws.Range("K105").Value = 3832.2856
$ws.Range("M105").Value = -2085.2856
$ws.Range("H134").Value = 5657.364
$ws.Range("I134").Value = 5657.364
$ws.Range("K134").Value = 16972.092
$ws.Range("M134").Value = -14437.092

# --- Sheet: CRP ---
$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H7").Value = 333385.12
$ws.Range("I7").Value = 51.545456
$ws.Range("J7").Value = 1250052.5
$ws.Range("K7").Value = 51.545456
$ws.Range("L7").Value = 1250052.5
$ws.Range("M7").Value = 61.454544
$ws.Range("N7").Value = -1250278.5
$ws.Range("H22").Value = 0
$ws.Range("I22").Value = 0
$ws.Range("K22").Value = 0
$ws.Range("M22").Value = ""
$ws.Range("H103").Value = 0
$ws.Range("I103").Value = 0
$ws.Range("K103").Value = 0
$ws.Range("M103").Value = ""

# --- Sheet: CUL ---
$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H4").Value = 2038.8966
$ws.Range("I4").Value = 2145.5833
$ws.Range("K4").Value = 6436.749899999999
$ws.Range("M4").Value = -6324.749899999999
$ws.Range("H7").Value = 289.66666
$ws.Range("I7").Value = 154
$ws.Range("J7").Value = 386.57144
$ws.Range("K7").Value = 462
$ws.Range("L7").Value = 1159.71432
$ws.Range("M7").Value = -350
$ws.Range("N7").Value = -1383.71432
$ws.Range("H11").Value = 222719.78
$ws.Range("I11").Value = 1000062
$ws.Range("J11").Value = 622
$ws.Range("K11").Value = 3000186
$ws.Range("L11").Value = 1866
$ws.Range("M11").Value = -3000046
$ws.Range("N11").Value = -2146
$ws.Range("H13").Value = 199
$ws.Range("J13").Value = 28.333334
$ws.Range("L13").Value = 85.00000199999999
$ws.Range("N13").Value = -421.000002
$ws.Range("H16").Value = 939
$ws.Range("I16").Value = 939
$ws.Range("J16").Value = 0
$ws.Range("K16").Value = 2817
$ws.Range("L16").Value = 0
$ws.Range("M16").Value = -2644
$ws.Range("N16").Value = ""
$ws.Range("H23").Value = 443.33334
$ws.Range("J23").Value = 285
$ws.Range("L23").Value = 855
$ws.Range("N23").Value = -1325
$ws.Range("H26").Value = 132068.12
$ws.Range("I26").Value = 201219.6
$ws.Range("K26").Value = 603658.8
$ws.Range("M26").Value = -603370.8
$ws.Range("H80").Value = 12979.583
$ws.Range("J80").Value = 13595.6
$ws.Range("L80").Value = 40786.8
$ws.Range("N80").Value = -42658.8
$ws.Range("H83").Value = 12979.583
$ws.Range("J83").Value = 13595.6
$ws.Range("L83").Value = 122360.4
$ws.Range("N83").Value = -131720.4

# --- Sheet: GSM ---
$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H2").Value = 151.375
$ws.Range("I2").Value = 123
$ws.Range("J2").Value = 350
$ws.Range("K2").Value = 123
$ws.Range("L2").Value = 350
$ws.Range("M2").Value = -10
$ws.Range("N2").Value = -576
$ws.Range("H3").Value = 300
$ws.Range("J3").Value = 300
$ws.Range("L3").Value = 300
$ws.Range("N3").Value = -532
$ws.Range("H97").Value = 854
$ws.Range("I97").Value = 842.7143
$ws.Range("J97").Value = 873.75
$ws.Range("K97").Value = 842.7143
$ws.Range("L97").Value = 873.75
$ws.Range("M97").Value = -346.7143
$ws.Range("N97").Value = -1865.75
$ws.Range("H122").Value = 2731.4
$ws.Range("I122").Value = 1436.25
$ws.Range("K122").Value = 4308.75
$ws.Range("M122").Value = -1858.75

# --- Sheet: LTW ---
$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H16").Value = 1398.3334
$ws.Range("J16").Value = 0
$ws.Range("L16").Value = 0
$ws.Range("N16").Value = ""
$ws.Range("H20").Value = 7857.143
$ws.Range("J20").Value = 8000
$ws.Range("L20").Value = 8000
$ws.Range("N20").Value = -8452
$ws.Range("H34").Value = 3250
$ws.Range("I34").Value = 3250
$ws.Range("K34").Value = 3250
$ws.Range("M34").Value = -3078
$ws.Range("H46").Value = 2462.2856
$ws.Range("I46").Value = 2077.6667
$ws.Range("J46").Value = 3154.6
$ws.Range("K46").Value = 2077.6667
$ws.Range("L46").Value = 3154.6
$ws.Range("M46").Value = -1889.6667
$ws.Range("N46").Value = -3530.6
$ws.Range("H61").Value = 9713.571
$ws.Range("J61").Value = 10001.75
$ws.Range("L61").Value = 10001.75
$ws.Range("N61").Value = -10405.75
$ws.Range("H82").Value = 1937.5
$ws.Range("J82").Value = 2875
$ws.Range("L82").Value = 2875
$ws.Range("N82").Value = -3597
$ws.Range("H85").Value = 1937.5
$ws.Range("J85").Value = 2875
$ws.Range("L85").Value = 2875
$ws.Range("N85").Value = -5371
$ws.Range("H113").Value = 9713.571
$ws.Range("J113").Value = 10001.75
$ws.Range("L113").Value = 10001.75
$ws.Range("N113").Value = -14341.75
$ws.Range("H122").Value = 3251.8333
$ws.Range("I122").Value = 3251.8333
$ws.Range("K122").Value = 9755.499899999999
$ws.Range("M122").Value = -7305.499899999999
$ws.Range("H132").Value = 3137.7778
$ws.Range("I132").Value = 2707.5
$ws.Range("K132").Value = 8122.5
$ws.Range("M132").Value = -5592.5

# --- Sheet: WVR ---
$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H96").Value = 4881.3335
$ws.Range("I96").Value = 5100
$ws.Range("K96").Value = 5100
$ws.Range("M96").Value = -3727
$ws.Range("H136").Value = 3739.743
$ws.Range("I136").Value = 4185.115
$ws.Range("J136").Value = 2453.111
$ws.Range("K136").Value = 12555.345
$ws.Range("L136").Value = 7359.333
$ws.Range("M136").Value = -10005.345
$ws.Range("N136").Value = -12459.333
